$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Replace the employee roster (rows 2-20). Columns A (ID), D (Age),
# E (Sexe) and H (Adresse) keep the same values as before; B (Nom),
# C (Prenom), F (Telephone) and G (Email) are replaced with the new
# dataset pulled from the Downloads copy of the datasource.
# ---------------------------------------------------------------------

$ws.Cells.Item(2,2).Value  = "Lorandelo"
$ws.Cells.Item(2,3).Value  = "Jordano"
$ws.Cells.Item(2,6).Value  = "778563565"

$ws.Cells.Item(3,2).Value  = "Derdourian"
$ws.Cells.Item(3,3).Value  = "Branahm"
$ws.Cells.Item(3,6).Value  = "771004541"

$ws.Cells.Item(4,2).Value  = "Trama"
$ws.Cells.Item(4,3).Value  = "mirella"
$ws.Cells.Item(4,6).Value  = "771003265"

$ws.Cells.Item(5,2).Value  = "Bigorno"
$ws.Cells.Item(5,3).Value  = "Jean-louis"
$ws.Cells.Item(5,6).Value  = "771013265"

$ws.Cells.Item(6,2).Value  = "Françoise"
$ws.Cells.Item(6,3).Value  = "Matthieu"
$ws.Cells.Item(6,6).Value  = "771003265"

$ws.Cells.Item(7,2).Value  = "Castor"
$ws.Cells.Item(7,3).Value  = "Jeanne"
$ws.Cells.Item(7,6).Value  = "771013265"

$ws.Cells.Item(8,2).Value  = "John"
$ws.Cells.Item(8,3).Value  = "Jean"
$ws.Cells.Item(8,6).Value  = "771003265"

$ws.Cells.Item(9,2).Value  = "Adjaye"
$ws.Cells.Item(9,3).Value  = "Brad"
$ws.Cells.Item(9,6).Value  = "771013265"

$ws.Cells.Item(10,2).Value = "Ansermina"
$ws.Cells.Item(10,3).Value = "Ivan"
$ws.Cells.Item(10,6).Value = "771013266"

$ws.Cells.Item(11,2).Value = "YANG"
$ws.Cells.Item(11,3).Value = "stecy"
$ws.Cells.Item(11,6).Value = "771013267"

$ws.Cells.Item(12,2).Value = "Briquello"
$ws.Cells.Item(12,3).Value = "Irenee"
$ws.Cells.Item(12,6).Value = "771013268"

$ws.Cells.Item(13,2).Value = "lafleur"
$ws.Cells.Item(13,3).Value = "henrina"
$ws.Cells.Item(13,6).Value = "771013269"

$ws.Cells.Item(14,2).Value = "BAYA"
$ws.Cells.Item(14,3).Value = "Annie"
$ws.Cells.Item(14,6).Value = "771013270"

$ws.Cells.Item(15,2).Value = "Diaz"
$ws.Cells.Item(15,3).Value = "julia"
$ws.Cells.Item(15,6).Value = "771013271"

$ws.Cells.Item(16,2).Value = "Kanto"
$ws.Cells.Item(16,3).Value = "Lafleur"
$ws.Cells.Item(16,6).Value = "771013272"

$ws.Cells.Item(17,2).Value = "Rodrigurez"
$ws.Cells.Item(17,3).Value = "Hernandez"
$ws.Cells.Item(17,6).Value = "771013273"

$ws.Cells.Item(18,2).Value = "Romano"
$ws.Cells.Item(18,3).Value = "Dimitrio"
$ws.Cells.Item(18,6).Value = "771013274"

$ws.Cells.Item(19,2).Value = "Rateau"
$ws.Cells.Item(19,3).Value = "Chloe"
$ws.Cells.Item(19,6).Value = "771013275"

$ws.Cells.Item(20,2).Value = "Rachida"
$ws.Cells.Item(20,3).Value = "Kadidja"
$ws.Cells.Item(20,6).Value = "771013276"

# ---------------------------------------------------------------------
# Rebuild column G (Email) + its mailto hyperlinks. Drop every existing
# hyperlink first so re-adding them doesn't leave stale duplicates
# behind, then re-create them in the same order the author did.
# ---------------------------------------------------------------------

$ws.Hyperlinks.Delete()

$ws.Range("G2").Value  = "ljordano@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:ljordano@gmail.com")

$ws.Range("G20").Value = "r.kadidja@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G20"), "mailto:r.kadidja@gmail.com")

$ws.Range("G19").Value = "r.chloe@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G19"), "mailto:r.chloe@gmail.com")

$ws.Range("G18").Value = "r.dimitrio@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G18"), "mailto:r.dimitrio@gmail.com")

$ws.Range("G17").Value = "r.hernandez@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G17"), "mailto:r.hernandez@gmail.com")

$ws.Range("G16").Value = "k.lafleur@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G16"), "mailto:k.lafleur@gmail.com")

$ws.Range("G15").Value = "d.julia@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G15"), "mailto:d.julia@gmail.com")

$ws.Range("G3").Value  = "d.branahm@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:d.branahm@gmail.com")

$ws.Range("G4").Value  = "t.mirella@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G4"), "mailto:t.mirella@gmail.com")

$ws.Range("G5").Value  = "b.jeanlouis@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G5"), "mailto:b.jeanlouis@gmail.com")

$ws.Range("G7").Value  = "cjeanne@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G7"), "mailto:cjeanne@gmail.com")

$ws.Range("G6").Value  = "fmatthieu@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G6"), "mailto:fmatthieu@gmail.com")

$ws.Range("G8").Value  = "jjean@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G8"), "mailto:jjean@gmail.com")

$ws.Range("G9").Value  = "abrad@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G9"), "mailto:abrad@gmail.com")

$ws.Range("G10").Value = "a.ivan@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G10"), "mailto:a.ivan@gmail.com")

$ws.Range("G11").Value = "yangstecy@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G11"), "mailto:yangstecy@gmail.com")

$ws.Range("G12").Value = "briquelloirenee@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G12"), "mailto:briquelloirenee@gmail.com")

$ws.Range("G13").Value = "lafleurhenrina@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G13"), "mailto:lafleurhenrina@gmail.com")

$ws.Range("G14").Value = "bannie@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G14"), "mailto:bannie@gmail.com")

# ---------------------------------------------------------------------
# Restore the current selection to the cell the author ended up on.
# ---------------------------------------------------------------------
$ws.Range("G13").Select()
